# Adds a "2022-Q1" sheet (fund holding detail) right after "2021-Q4", and
# refreshes the "总计" (summary) sheet with a new row for 2022-Q1, matching
# the structure/style already used by the existing sheets.

$wb = $excel.ActiveWorkbook

$wsQ4   = $wb.Worksheets.Item(1)   # "2021-Q4" - stays as-is
$wsOldTotal = $wb.Worksheets.Item(2)   # old "总计" sheet - will be rebuilt

# Remember the old summary rows before the sheet is removed & recreated.
$oldDate  = $wsOldTotal.Range("B2").Value2
$oldCount = $wsOldTotal.Range("C2").Value2
$oldValue = $wsOldTotal.Range("D2").Value2

$wsOldTotal.Delete() | Out-Null

# --- New sheet: 2022-Q1 (fund holding detail), placed right after 2021-Q4 ---
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

$q1Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q1Header.Length; $i++) {
    $cell = $wsQ1.Cells.Item(1, 2 + $i)
    $cell.Value2 = $q1Header[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Row 2 - 005482
$wsQ1.Cells.Item(2, 1).Value2 = 0
$wsQ1.Cells.Item(2, 1).Font.Bold = $true
$wsQ1.Cells.Item(2, 1).Borders.LineStyle = 1
$wsQ1.Cells.Item(2, 1).HorizontalAlignment = -4108
$wsQ1.Cells.Item(2, 1).VerticalAlignment = -4160

$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("B2").Value2 = "005482"
$wsQ1.Range("C2").Value2 = "博时创新驱动灵活配置混合A"
$wsQ1.Range("D2").NumberFormat = "@"
$wsQ1.Range("D2").Value2 = "0.41"
$wsQ1.Range("E2").NumberFormat = "@"
$wsQ1.Range("E2").Value2 = "88.65"
$wsQ1.Range("F2").NumberFormat = "@"
$wsQ1.Range("F2").Value2 = "8.90"
$wsQ1.Range("G2").NumberFormat = "@"
$wsQ1.Range("G2").Value2 = "0.0365"
$wsQ1.Range("H2").Value2 = 1

# Row 3 - 005483
$wsQ1.Cells.Item(3, 1).Value2 = 1
$wsQ1.Cells.Item(3, 1).Font.Bold = $true
$wsQ1.Cells.Item(3, 1).Borders.LineStyle = 1
$wsQ1.Cells.Item(3, 1).HorizontalAlignment = -4108
$wsQ1.Cells.Item(3, 1).VerticalAlignment = -4160

$wsQ1.Range("B3").NumberFormat = "@"
$wsQ1.Range("B3").Value2 = "005483"
$wsQ1.Range("C3").Value2 = "博时创新驱动灵活配置混合C"
$wsQ1.Range("D3").NumberFormat = "@"
$wsQ1.Range("D3").Value2 = "0.05"
$wsQ1.Range("E3").NumberFormat = "@"
$wsQ1.Range("E3").Value2 = "88.65"
$wsQ1.Range("F3").NumberFormat = "@"
$wsQ1.Range("F3").Value2 = "8.90"
$wsQ1.Range("G3").NumberFormat = "@"
$wsQ1.Range("G3").Value2 = "0.0044"
$wsQ1.Range("H3").Value2 = 1

# --- New sheet: 总计 (summary), placed right after 2022-Q1 ---
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

$totalHeader = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeader.Length; $i++) {
    $cell = $wsTotal.Cells.Item(1, 2 + $i)
    $cell.Value2 = $totalHeader[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Row 2 - newest quarter (2022-Q1) goes first
$wsTotal.Cells.Item(2, 1).Value2 = 0
$wsTotal.Cells.Item(2, 1).Font.Bold = $true
$wsTotal.Cells.Item(2, 1).Borders.LineStyle = 1
$wsTotal.Cells.Item(2, 1).HorizontalAlignment = -4108
$wsTotal.Cells.Item(2, 1).VerticalAlignment = -4160

$wsTotal.Range("B2").Value2 = "2022-Q1"
$wsTotal.Range("C2").Value2 = 2
$wsTotal.Range("D2").Value2 = 0.04

# Row 3 - previous quarter (2021-Q4), carried over from the old summary sheet
$wsTotal.Cells.Item(3, 1).Value2 = 1
$wsTotal.Cells.Item(3, 1).Font.Bold = $true
$wsTotal.Cells.Item(3, 1).Borders.LineStyle = 1
$wsTotal.Cells.Item(3, 1).HorizontalAlignment = -4108
$wsTotal.Cells.Item(3, 1).VerticalAlignment = -4160

$wsTotal.Range("B3").Value2 = $oldDate
$wsTotal.Range("C3").Value2 = $oldCount
$wsTotal.Range("D3").Value2 = $oldValue
